$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet1 ("Minimum"): rebuild the header row with the new field set.
#    Old:  A..H = mediaAsset0, mediaAsset0.keywords, mediaAsset1,
#                 mediaAsset1.keywords, individualID, genus, specificEpithet,
#                 submitterID
#    New:  A..L = mediaAsset0, mediaAsset0.keywords, mediaAsset1,
#                 mediaAsset1.keywords, individualID, locationID, genus,
#                 specificEpithet, year, month, day, submitterID
#
# Work right-to-left / capture-formatting-before-overwrite so no cell's
# original look gets clobbered before it has been copied onward.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Minimum")

# submitterID moves from H to L - grab H's current (distinct) look first.
$ws1.Range("H1").Copy()
$ws1.Range("L1").PasteSpecial(-4122)
$ws1.Range("L1").Value = "Encounter.submitterID"

# specificEpithet moves from G to H; it should pick up the shared
# "highlighted header" look that already sits on F/G, not H's old look.
$ws1.Range("G1").Copy()
$ws1.Range("H1").PasteSpecial(-4122)
$ws1.Range("H1").Value = "Encounter.specificEpithet"

# genus moves from F to G - G already carries the right look, just retarget.
$ws1.Range("G1").Value = "Encounter.genus"

# locationID is brand new at F - F already carries the right look.
$ws1.Range("F1").Value = "Encounter.locationID"

# year is brand new at I, reusing the shared highlighted-header look.
$ws1.Range("F1").Copy()
$ws1.Range("I1").PasteSpecial(-4122)
$ws1.Range("I1").Value = "Encounter.year"

# month/day are brand new at J/K with a text number format plus a white fill.
$ws1.Range("J1").NumberFormat = "@"
$ws1.Range("J1").Interior.ThemeColor = 2
$ws1.Range("J1").Copy()
$ws1.Range("K1").PasteSpecial(-4122)
$ws1.Range("J1").Value = "Encounter.month"
$ws1.Range("K1").Value = "Encounter.day"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Rename the second sheet and replace its content with the new "help
#    links" layout.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Explanations")
$ws2.Name = "Info"

# Wipe the old explanatory table entirely.
$ws2.Cells.Clear()

$ws2.Range("A1").Value = "For a list of available fields and how to use them see:"
$ws2.Range("A2").Value = "Find demos and tutorials on our YouTube channel: "

$ws2.Hyperlinks.Add($ws2.Range("B1"), "https://wildbook.docs.wildme.org/data/bulk-import-beta.html", "fields-available")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://www.youtube.com/@wildme3451/videos")

$excel.CutCopyMode = $false

Write-Host "edit applied"
